$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 15364.054
$ws.Range("J70").Value = 17105.576
$ws.Range("L70").Value = 51316.728
$ws.Range("N70").Value = -51856.728
$ws.Range("H73").Value = 15364.054
$ws.Range("J73").Value = 17105.576
$ws.Range("L73").Value = 51316.728
$ws.Range("N73").Value = -53188.728
$ws.Range("H113").Value = 3847824
$ws.Range("I113").Value = 5883936
$ws.Range("J113").Value = 1834.6666
$ws.Range("K113").Value = 5883936
$ws.Range("L113").Value = 1834.6666
$ws.Range("M113").Value = -5880682
$ws.Range("N113").Value = -8342.6666
$ws.Range("H116").Value = 29178030
$ws.Range("I116").Value = 16669141
$ws.Range("J116").Value = 41686920
$ws.Range("K116").Value = 16669141
$ws.Range("L116").Value = 41686920
$ws.Range("M116").Value = -16665699
$ws.Range("N116").Value = -41693804
$ws.Range("H132").Value = 5292479
$ws.Range("I132").Value = 1547.55
$ws.Range("K132").Value = 4642.65
$ws.Range("M132").Value = -2112.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7535898.5
$ws.Range("I32").Value = 1897875.6
$ws.Range("K32").Value = 1897875.6
$ws.Range("M32").Value = -1897588.6
$ws.Range("H45").Value = 455731.72
$ws.Range("I45").Value = 910072.8
$ws.Range("K45").Value = 910072.8
$ws.Range("M45").Value = -909695.8
$ws.Range("H61").Value = 2139011.5
$ws.Range("J61").Value = 7353713
$ws.Range("L61").Value = 7353713
$ws.Range("N61").Value = -7354137
$ws.Range("H74").Value = 51587972
$ws.Range("I74").Value = 42262296
$ws.Range("J74").Value = 88890670
$ws.Range("K74").Value = 42262296
$ws.Range("L74").Value = 88890670
$ws.Range("M74").Value = -42261422
$ws.Range("N74").Value = -88892418
$ws.Range("H77").Value = 51587972
$ws.Range("I77").Value = 42262296
$ws.Range("J77").Value = 88890670
$ws.Range("K77").Value = 211311480
$ws.Range("L77").Value = 444453350
$ws.Range("M77").Value = -211307112
$ws.Range("N77").Value = -444462086
$ws.Range("H122").Value = 4666.6665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4666.6665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13999.9995
$ws.Range("N122").Value = -18899.9995
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 11114942
$ws.Range("I132").Value = 11115191
$ws.Range("J132").Value = 11113451
$ws.Range("K132").Value = 33345573
$ws.Range("L132").Value = 33340353
$ws.Range("M132").Value = -33343043
$ws.Range("N132").Value = -33345413
$ws.Range("H136").Value = 2139011.5
$ws.Range("J136").Value = 7353713
$ws.Range("L136").Value = 22061139
$ws.Range("N136").Value = -22066239

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 961.1111
$ws.Range("I107").Value = 847.25
$ws.Range("J107").Value = 993.6429000000001
$ws.Range("K107").Value = 847.25
$ws.Range("L107").Value = 993.6429000000001
$ws.Range("M107").Value = 1072.75
$ws.Range("N107").Value = -4833.6429
$ws.Range("H134").Value = 9741616
$ws.Range("I134").Value = 12196155
$ws.Range("J134").Value = 2553321.2
$ws.Range("K134").Value = 36588465
$ws.Range("L134").Value = 7659963.600000001
$ws.Range("M134").Value = -36585930
$ws.Range("N134").Value = -7665033.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2267545.8
$ws.Range("I31").Value = 1264067.2
$ws.Range("J31").Value = 4814837
$ws.Range("K31").Value = 1264067.2
$ws.Range("L31").Value = 4814837
$ws.Range("M31").Value = -1263772.2
$ws.Range("N31").Value = -4815427
$ws.Range("H34").Value = 2267545.8
$ws.Range("I34").Value = 1264067.2
$ws.Range("J34").Value = 4814837
$ws.Range("K34").Value = 1264067.2
$ws.Range("L34").Value = 4814837
$ws.Range("M34").Value = -1263865.2
$ws.Range("N34").Value = -4815241
$ws.Range("H94").Value = 125009870
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 125009870
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 125009870
$ws.Range("N94").Value = -125010772
$ws.Range("M94").ClearContents()
$ws.Range("H107").Value = 743
$ws.Range("I107").Value = 330
$ws.Range("J107").Value = 855.63635
$ws.Range("K107").Value = 330
$ws.Range("L107").Value = 855.63635
$ws.Range("M107").Value = 1590
$ws.Range("N107").Value = -4695.63635
$ws.Range("H122").Value = 5204.5884
$ws.Range("I122").Value = 10218.75
$ws.Range("J122").Value = 747.55554
$ws.Range("K122").Value = 30656.25
$ws.Range("L122").Value = 2242.66662
$ws.Range("M122").Value = -28206.25
$ws.Range("N122").Value = -7142.66662
$ws.Range("H134").Value = 2868057.2
$ws.Range("I134").Value = 18404
$ws.Range("J134").Value = 5005297
$ws.Range("K134").Value = 55212
$ws.Range("L134").Value = 15015891
$ws.Range("M134").Value = -52677
$ws.Range("N134").Value = -15020961

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 870.18463
$ws.Range("J131").Value = 927.0862
$ws.Range("L131").Value = 2781.2586
$ws.Range("N131").Value = -12861.2586
$ws.Range("H132").Value = 2618.2727
$ws.Range("I132").Value = 2901.6
$ws.Range("J132").Value = 2382.1667
$ws.Range("K132").Value = 26114.4
$ws.Range("L132").Value = 21439.5003
$ws.Range("M132").Value = -23584.4
$ws.Range("N132").Value = -26499.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8269.429
$ws.Range("I102").Value = 9272.666999999999
$ws.Range("K102").Value = 9272.666999999999
$ws.Range("M102").Value = -7650.666999999999
$ws.Range("H107").Value = 161.91667
$ws.Range("I107").Value = 82.55556
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 82.55556
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1837.44444
$ws.Range("N107").Value = -4240
$ws.Range("H122").Value = 9836982
$ws.Range("I122").Value = 42365.152
$ws.Range("K122").Value = 127095.456
$ws.Range("M122").Value = -124645.456
$ws.Range("H132").Value = 26043260
$ws.Range("I132").Value = 61905290
$ws.Range("J132").Value = 10104581
$ws.Range("K132").Value = 185715870
$ws.Range("L132").Value = 30313743
$ws.Range("M132").Value = -185713340
$ws.Range("N132").Value = -30318803

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1107.72
$ws.Range("I7").Value = 873.13336
$ws.Range("K7").Value = 873.13336
$ws.Range("M7").Value = -761.13336
$ws.Range("H40").Value = 2369.5334
$ws.Range("I40").Value = 1296.5
$ws.Range("K40").Value = 1296.5
$ws.Range("M40").Value = -1160.5
$ws.Range("H122").Value = 6510476
$ws.Range("I122").Value = 888716.7
$ws.Range("J122").Value = 20002698
$ws.Range("K122").Value = 2666150.1
$ws.Range("L122").Value = 60008094
$ws.Range("M122").Value = -2663700.1
$ws.Range("N122").Value = -60012994
$ws.Range("H126").Value = 1107.72
$ws.Range("I126").Value = 873.13336
$ws.Range("K126").Value = 2619.40008
$ws.Range("M126").Value = -149.4000800000003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1246.909
$ws.Range("I122").Value = 1167.238
$ws.Range("J122").Value = 1386.3334
$ws.Range("K122").Value = 3501.714
$ws.Range("L122").Value = 4159.0002
$ws.Range("M122").Value = -1051.714
$ws.Range("N122").Value = -9059.0002
